$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 397, shifting existing rows 397:485 down to 398:486
$ws.Rows(397).Insert()

# Populate the newly inserted row 397 with the new data record
$ws.Range("A397").Value = 9
$ws.Range("B397").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C397").Value = "Metropolitana"
$ws.Range("D397").Value = 44889
$ws.Range("D397").NumberFormat = $ws.Range("D398").NumberFormat
$ws.Range("E397").Value = 13
$ws.Range("F397").Value = 100112039
$ws.Range("G397").Value = "Ciboulette"
$ws.Range("H397").Value = "Sin especificar"
$ws.Range("I397").Value = "Primera"
$ws.Range("J397").Value = 440
$ws.Range("K397").Value = 800
$ws.Range("L397").Value = 1000
$ws.Range("M397").Value = 891
$ws.Range("N397").Value = "`$/docena de atados"
$ws.Range("O397").Value = "Región Metropolitana"
$ws.Range("P397").Value = 297
$ws.Range("Q397").Value = 3
$ws.Range("R397").Value = "Hortaliza"
